$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Word's "spell-check" left stray <w:proofErr> markers splitting some
# cells' text across multiple runs. The cleanest way to get Word to
# regenerate a cell's paragraph/run markup from scratch (dropping the
# proofErr wrapper entirely) is to insert a brand-new row with the
# desired text and delete the old one in its place.
function Rebuild-Row($table, $rowIndex, $c1, $c2, $c3) {
    $oldRow = $table.Rows.Item($rowIndex)
    $beforeRow = $table.Rows.Item($rowIndex + 1)
    $newRow = $table.Rows.Add($beforeRow)
    $newRow.Cells.Item(1).Range.Text = $c1
    $newRow.Cells.Item(2).Range.Text = $c2
    $newRow.Cells.Item(3).Range.Text = $c3
    $oldRow.Delete()
}

# Rows whose "Kuvaus/Selite" cell had spurious proofErr-split runs;
# rebuilt here so the text collapses back into a single <w:r>.
Rebuild-Row $t 4  "03.02.2023" "0.5" "Weekly"
Rebuild-Row $t 5  "09.02.2023" "2"   "Class diagram ja tehtävänluonnin UI"
Rebuild-Row $t 6  "10.02.2023" "0.5" "Weekly"
Rebuild-Row $t 8  "17.02.2023" "3"   "Weekly + uuden tagin luominen"
Rebuild-Row $t 9  "24.02.2023" "1"   "Weekly + tekijän muutoksia "
Rebuild-Row $t 10 "28.02.2023" "3"   "Uuden tehtävän luonnin muutoksia: tägit ja tekijä "
Rebuild-Row $t 11 "03.03.2023" "2"   "Weekly + tägien muutoksia "

# Fill in the next two previously-blank rows with the new entries for
# the final two sessions ("loppupaltsun kahden osion slaidit lisätty").
$t.Cell(12, 1).Range.Text = "06.03.2023"
$t.Cell(12, 2).Range.Text = "1"
$t.Cell(12, 3).Range.Text = "Kaavioiden hiomista"

$t.Cell(13, 1).Range.Text = "08.03.2023"
$t.Cell(13, 2).Range.Text = "1"
$t.Cell(13, 3).Range.Text = "Lopputapaamiseen valmistautumista"
